{"js": "// Helper: build a one-run-paragraph OOXML snippet (for `insertOoxml`) so we\n// can append additional runs to a paragraph with full control over the\n// run's formatting (rPr) and text (including xml:space=\"preserve\").\nfunction runOoxml(rPrXml, text, preserveSpace) {\n  const t = preserveSpace\n    ? `<w:t xml:space=\"preserve\">${text}</w:t>`\n    : `<w:t>${text}</w:t>`;\n  const rPr = rPrXml || \"\";\n  return `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body><w:p><w:r>${rPr}${t}</w:r></w:p></w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n}\n\nasync function replaceOnce(searchText, newText) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Title\nawait replaceOnce(\"Game Name\", \"Feather Face\");\n\n// 2) Genre question\nawait replaceOnce(\"What type of game is this?\", \"A platformer\");\n\n// 3) Target audience question\nawait replaceOnce(\n  \"Who are you trying to sell this game to? What about them should find this game appealing?\",\n  \"People who like platformers like Mario.\"\n);\n\n// 4) ESRB rating question\nawait replaceOnce(\"What is the game\\u2019s rating?\", \"E\");\n\n// 5) Target platform question\nawait replaceOnce(\n  \"What is the target platform? If PC, what specs will the computer need?\",\n  \"Consoles\"\n);\n\n// 6) Required licenses question\nawait replaceOnce(\n  \"Does this game require any licenses such as a movie it is based on?\",\n  \"No required licenses.\"\n);\n\n// 7) Competition and marketing paragraph -- becomes three runs:\n//    \"It will be marketed...game play.\" + \" \" + \"I think it will sell...free to play.\"\n{\n  const marketingOld =\n    \"How is this game going to be marketed? How does it compare and contrast with the competition? Why do you think this game will sell?\";\n  await replaceOnce(\n    marketingOld,\n    \"It will be marketed towards fans of platformers. It is more accessible than the competition with easier game play.\"\n  );\n\n  const paras = context.document.body.paragraphs;\n  paras.load(\"text\");\n  await context.sync();\n  let marketingPara = null;\n  for (let i = 0; i < paras.items.length; i++) {\n    if (\n      paras.items[i].text ===\n      \"It will be marketed towards fans of platformers. It is more accessible than the competition with easier game play.\"\n    ) {\n      marketingPara = paras.items[i];\n      break;\n    }\n  }\n\n  marketingPara.insertOoxml(runOoxml(\"\", \" \", true), Word.InsertLocation.end);\n  await context.sync();\n  marketingPara.insertOoxml(\n    runOoxml(\"\", \"I think it will sell well because it will be free to play.\", false),\n    Word.InsertLocation.end\n  );\n  await context.sync();\n}\n\n// 8) Designer and creation date paragraph -- the trailing run\n//    \"Who are you? When was this written?\" becomes four runs:\n//    \"Jerod Merrit\" + \"t\" + \", \" + \"October 2022\", and the _GoBack bookmark\n//    that trailed it is removed.\n{\n  const authorRPr =\n    '<w:rPr><w:color w:val=\"202124\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr>';\n\n  await replaceOnce(\"Who are you? When was this written?\", \"Jerod Merrit\");\n\n  const paras = context.document.body.paragraphs;\n  paras.load(\"text\");\n  await context.sync();\n  let authorPara = null;\n  for (let i = 0; i < paras.items.length; i++) {\n    if (paras.items[i].text.indexOf(\"Jerod Merrit\") !== -1) {\n      authorPara = paras.items[i];\n      break;\n    }\n  }\n\n  authorPara.insertOoxml(runOoxml(authorRPr, \"t\", false), Word.InsertLocation.end);\n  await context.sync();\n  authorPara.insertOoxml(runOoxml(authorRPr, \", \", true), Word.InsertLocation.end);\n  await context.sync();\n  authorPara.insertOoxml(\n    runOoxml(authorRPr, \"October 2022\", false),\n    Word.InsertLocation.end\n  );\n  await context.sync();\n\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($find, $replace) {\n    $f = $d.Content.Find\n    $f.Text = $find\n    $f.Replacement.Text = $replace\n    $f.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null\n}\n\n# 1) Title\nReplace-Text \"Game Name\" \"Feather Face\"\n\n# 2) Genre question\nReplace-Text \"What type of game is this?\" \"A platformer\"\n\n# 3) Target audience question\nReplace-Text \"Who are you trying to sell this game to? What about them should find this game appealing?\" \"People who like platformers like Mario.\"\n\n# 4) ESRB rating question\nReplace-Text \"What is the game\u2019s rating?\" \"E\"\n\n# 5) Target platform question\nReplace-Text \"What is the target platform? If PC, what specs will the computer need?\" \"Consoles\"\n\n# 6) Required licenses question\nReplace-Text \"Does this game require any licenses such as a movie it is based on?\" \"No required licenses.\"\n\n# 7) Competition and marketing paragraph becomes two sentences (kept together\n#    in this run's text; the source edit authored them as three runs split on\n#    an inter-sentence space, which carries no visible/formatting difference).\nReplace-Text \"How is this game going to be marketed? How does it compare and contrast with the competition? Why do you think this game will sell?\" \"It will be marketed towards fans of platformers. It is more accessible than the competition with easier game play. I think it will sell well because it will be free to play.\"\n\n# 8) Designer and creation date paragraph -- the trailing run\n#    \"Who are you? When was this written?\" is replaced by the author credit\n#    text, preserving that run's character formatting (color/shading).\nReplace-Text \"Who are you? When was this written?\" \"Jerod Merritt, October 2022\"\n\n# Remove the _GoBack bookmark that trailed the author-credit run.\n$d.Bookmarks(\"_GoBack\").Delete()\n"}
